# Update "Tiempos x Recurso x Sprint" workbook:
# Refresh the estimated-SP figures (and the dependent totals) for Sprint 2/3/4
# in both the summary table (B1:H9 / Tabla1) and the detail grid (A19:N27),
# then leave the selection where the author last clicked (H12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Summary table (Tabla1, B1:H9) -----------------------------------------
# Column D = "Estimated SP Sprint 2"
$ws.Range("D2").Value = 41
$ws.Range("D3").Value = 22
$ws.Range("D4").Value = 55
$ws.Range("D5").Value = 72
$ws.Range("D6").Value = 28

# Column E = "Estimated SP Sprint 3"
$ws.Range("E4").Value = 26
$ws.Range("E5").Value = 48

# Column F = "Estimated SP Sprint 4"
$ws.Range("F5").Value = 32

# --- Detail grid (A19:N27) mirrors the same figures -------------------------
# Column E mirrors Tabla1 column D
$ws.Range("E20").Value = 41
$ws.Range("E21").Value = 22
$ws.Range("E22").Value = 55
$ws.Range("E23").Value = 72
$ws.Range("E24").Value = 28

# Column G mirrors Tabla1 column E
$ws.Range("G22").Value = 26
$ws.Range("G23").Value = 48

# Column I mirrors Tabla1 column F
$ws.Range("I23").Value = 32

# Totals (H column / M column) and the SUM rows (9 / 27) are formulas and
# recalculate automatically.

# --- Leave the selection where it was when the file was last saved ---------
$ws.Range("H12").Select() | Out-Null
